# Refresh the crypto price/volume table (GitHub Actions bot update).
#
# Price (column D) and Volume(1h) (column E) cells in this sheet are stored
# as plain text (e.g. "25.985.73", "  +0.76%  "), not numbers. Several D
# values parse as plain floats (e.g. "246.01"), so a bare
# `Range.Value = "246.01"` would silently coerce them to a Double and lose
# the original text formatting. Prefixing with an apostrophe forces Excel
# to keep them as text, and resetting `.Style` back to "Normal" afterwards
# strips the "number stored as text" quote-prefix styling that trick adds,
# so the cell ends up with no style index at all - matching the source
# cells exactly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.985.73"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.76%  "
$ws.Range("D3").Value = "'1.740.44"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'246.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.59%  "
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("D7").Value = "'0.5024"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.33%  "
$ws.Range("D8").Value = "'0.2732"
$ws.Range("D8").Style = "Normal"
$ws.Range("E9").Value = "  +1.65%  "
$ws.Range("D10").Value = "'1.742.35"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.49%  "
$ws.Range("D11").Value = "'0.07256"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.71%  "
$ws.Range("D12").Value = "'0.6532"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.08%  "
$ws.Range("D13").Value = "'15.07"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.51%  "
$ws.Range("D14").Value = "'4.676"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.41%  "
$ws.Range("D15").Value = "'77.50"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.96%  "
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("D18").Value = "'26.013.14"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.79%  "
$ws.Range("D19").Value = "'11.91"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.27%  "
$ws.Range("D20").Value = "'0.000006854"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.90%  "
$ws.Range("D21").Value = "'1.968.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.40%  "
$ws.Range("D22").Value = "'4.493"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.22%  "
$ws.Range("D23").Value = "'8.701"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.59%  "
$ws.Range("D24").Value = "'5.389"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.74%  "
$ws.Range("D25").Value = "'135.54"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.58%  "
$ws.Range("D26").Value = "'1.509"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.89%  "
$ws.Range("D27").Value = "'15.23"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.88%  "
$ws.Range("D28").Value = "'1.780"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.82%  "
$ws.Range("D29").Value = "'105.36"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.58%  "
$ws.Range("D30").Value = "'3.947"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.00%  "
$ws.Range("D31").Value = "'0.08147"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.96%  "
$ws.Range("D32").Value = "'3.668"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.22%  "
$ws.Range("D33").Value = "'0.04701"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.90%  "
$ws.Range("D34").Value = "'2.668"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.38%  "
$ws.Range("D35").Value = "'0.9934"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.00%  "
$ws.Range("E36").Value = "  -0.86%  "
$ws.Range("E37").Value = "  +3.02%  "
$ws.Range("E38").Value = "  +2.13%  "
$ws.Range("D39").Value = "'1.922"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.74%  "
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("D41").Value = "'101.05"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.01%  "
$ws.Range("D42").Value = "'0.7912"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +8.41%  "
$ws.Range("D43").Value = "'0.3896"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.23%  "
$ws.Range("D44").Value = "'5.009"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.06%  "
$ws.Range("D45").Value = "'0.1165"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.30%  "
$ws.Range("D46").Value = "'6.309"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.91%  "
$ws.Range("D47").Value = "'55.64"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.62%  "
$ws.Range("E48").Value = "  +0.49%  "
$ws.Range("D49").Value = "'30.79"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.66%  "
$ws.Range("D50").Value = "'7.638"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.33%  "
$ws.Range("E51").Value = "  +2.20%  "
